{"js": "// The document contains several \"<id>...</id>\" tag pairs, each split\n// across three runs with different formatting:\n//   run1 \"<id>\"        (Courier New, color 7f6000, sz 18)\n//   run2 \"p160v_aN\"     (plain, color 000000)\n//   run3 \"</id>\"        (Courier New, color 7f6000, sz 18)\n// The edit collapses each triple into a single run (keeping run1's\n// Courier-New formatting) whose text drops the stray \"a\" before the\n// trailing number, e.g. \"<id>p160v_a1</id>\" -> \"<id>p160v_1</id>\".\n\nconst body = context.document.body;\n\nconst starts = body.search(\"<id>\", { matchCase: true, matchWildcards: false });\nconst ends = body.search(\"</id>\", { matchCase: true, matchWildcards: false });\nstarts.load(\"text\");\nends.load(\"text\");\nawait context.sync();\n\nconst count = Math.min(starts.items.length, ends.items.length);\n\n// Pair up each \"<id>\" with the \"</id>\" that closes it and expand a\n// range that spans all the runs in between (inclusive).\nconst fullRanges = [];\nfor (let i = 0; i < count; i++) {\n  fullRanges.push(starts.items[i].expandTo(ends.items[i]));\n}\nfullRanges.forEach((r) => r.load(\"text\"));\nawait context.sync();\n\n// Rewrite each \"<id>...</id>\" range as a single run of text, stripping\n// a literal \"a\" that immediately precedes the trailing digits of the id\n// (p160v_a1 -> p160v_1, p160v_a2 -> p160v_2, ...).\nfor (const range of fullRanges) {\n  const match = /^<id>([\\s\\S]*)<\\/id>$/.exec(range.text);\n  if (!match) continue;\n  const oldInner = match[1];\n  const newInner = oldInner.replace(/a(\\d+)$/, \"$1\");\n  if (newInner === oldInner) continue;\n  range.insertText(`<id>${newInner}</id>`, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document contains several \"<id>...</id>\" tag pairs, each split\n# across three runs with different formatting:\n#   run1 \"<id>\"        (Courier New, color 7f6000, sz 18)\n#   run2 \"p160v_aN\"     (plain, color 000000)\n#   run3 \"</id>\"        (Courier New, color 7f6000, sz 18)\n# The edit collapses each triple into a single run (keeping run1's\n# Courier-New formatting) whose text drops the stray \"a\" before the\n# trailing number, e.g. \"<id>p160v_a1</id>\" -> \"<id>p160v_1</id>\".\n\n$d = $word.ActiveDocument\n\n$searchStart = 0\n$storyEnd = $d.Content.End\n$count = 0\n\nwhile ($true) {\n  # Locate the next \"<id>\" opening tag.\n  $openRng = $d.Range($searchStart, $storyEnd)\n  $openRng.Find.ClearFormatting()\n  $openRng.Find.Text = \"<id>\"\n  $openRng.Find.Forward = $true\n  $openRng.Find.MatchWildcards = $false\n  $openRng.Find.Wrap = 0\n  $foundOpen = $openRng.Find.Execute()\n  if (-not $foundOpen) { break }\n\n  $idStart = $openRng.Start\n\n  # Locate the \"</id>\" closing tag that follows it.\n  $closeRng = $d.Range($openRng.End, $storyEnd)\n  $closeRng.Find.ClearFormatting()\n  $closeRng.Find.Text = \"</id>\"\n  $closeRng.Find.Forward = $true\n  $closeRng.Find.MatchWildcards = $false\n  $closeRng.Find.Wrap = 0\n  $foundClose = $closeRng.Find.Execute()\n  if (-not $foundClose) { break }\n\n  # Range spanning \"<id>...</id>\" across all the runs in between.\n  $fullRng = $d.Range($idStart, $closeRng.End)\n  $oldText = $fullRng.Text\n\n  if ($oldText -match '^<id>([\\s\\S]*)a(\\d+)</id>$') {\n    $newText = \"<id>\" + $matches[1] + $matches[2] + \"</id>\"\n    $fullRng.Text = $newText\n    $count = $count + 1\n    $searchStart = $idStart + $newText.Length\n  } else {\n    $searchStart = $closeRng.End\n  }\n\n  $storyEnd = $d.Content.End\n}\n\nWrite-Output \"replacements=$count\"\n"}
